# Fix the JSON-like parameter strings so they are valid JSON (quoted keys/values)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = '{"username":"SwatiChetty","password":"123456"}'
$ws.Range("D2").Value = '{"textToValidate":"Welcome Vivek!!"}'
$ws.Range("C3").Value = '{"username":"Swati","password":"123"}'
$ws.Range("D3").Value = '{"textToValidate":"Order Create Successfully"}'

# Widen column C now that its content is longer (closest attainable width to 46.453125)
$ws.Columns.Item(3).ColumnWidth = 45.65

# Leave the cursor on C9, matching the saved selection in the workbook
$ws.Range("C9").Select()
